$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AAT_Roles")

# New rows of NG Field / AAT Link / AAT Label data appended under the
# existing table (rows 2-11), continuing at row 12. Rows 12-15 were typed
# straight across (A, B, C per row). Rows 16-18 cover "Art historian",
# "Prof." and "Prof. Dr." - all three job-title labels (column A) were
# entered first, then the AAT link/label pair for row 16, then the pair
# for row 17, and row 18 simply reuses the same AAT link/label as row 17.
$ws.Cells.Item(12, 1).Value = "Scientist"
$ws.Cells.Item(12, 2).Value = "aat:300025788"
$ws.Cells.Item(12, 3).Value = "scientists"

$ws.Cells.Item(13, 1).Value = "Intern"
$ws.Cells.Item(13, 2).Value = "aat:300025902"
$ws.Cells.Item(13, 3).Value = "interns"

$ws.Cells.Item(14, 1).Value = "Chemical-technical assistant"
$ws.Cells.Item(14, 2).Value = "aat:300025898"
$ws.Cells.Item(14, 3).Value = "assistants"

$ws.Cells.Item(15, 1).Value = "Conservator"
$ws.Cells.Item(15, 2).Value = "aat:300102842"
$ws.Cells.Item(15, 3).Value = "conservators (people in conservation)"

$ws.Cells.Item(16, 1).Value = "Art historian"
$ws.Cells.Item(17, 1).Value = "Prof."
$ws.Cells.Item(18, 1).Value = "Prof. Dr."

$ws.Cells.Item(16, 2).Value = "aat:300025541"
$ws.Cells.Item(16, 3).Value = "art historians"

$ws.Cells.Item(17, 2).Value = "aat:300025533"
$ws.Cells.Item(17, 3).Value = "professors (teachers)"

$ws.Cells.Item(18, 2).Value = "aat:300025533"
$ws.Cells.Item(18, 3).Value = "professors (teachers)"

# Widen column A to fit the newly added, longer labels (best-fit to the
# longest new entry, "Chemical-technical assistant").
$ws.Columns.Item(1).ColumnWidth = 23.5

# Match the saved selection from the authoring session.
$ws.Range("F17").Select()
